# TORCH narrative + figures
#
# 1) Bump the cached "datetimeFigureOut" field text (24/06/2020 -> 27/06/2020)
#    on the slide master and every slide layout's date placeholder.
# 2) Shrink/relabel the two "includes (reflexes)" callouts on the only
#    slide down to just "reflexes", nudging the first one to the right.

$p = $ppt.ActivePresentation

# --- 1) Date placeholders on master + all custom layouts ---------------
$master = $p.SlideMaster

for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        if ($shp.TextFrame.HasText) {
            if ($shp.TextFrame.TextRange.Text -eq "24/06/2020") {
                $shp.TextFrame.TextRange.Text = "27/06/2020"
            }
        }
    }
}

for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        $shp = $layout.Shapes.Item($j)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.HasText) {
                if ($shp.TextFrame.TextRange.Text -eq "24/06/2020") {
                    $shp.TextFrame.TextRange.Text = "27/06/2020"
                }
            }
        }
    }
}

# --- 2) The two "includes (reflexes)" textboxes on slide 1 --------------
$s = $p.Slides.Item(1)

# "ZoneTexte 126" - moves right and shrinks, text shortened to "reflexes"
$sh17 = $s.Shapes.Item("ZoneTexte 126")
$sh17.TextFrame.TextRange.Text = "reflexes"
$sh17.Left = 347.6996062992126
$sh17.Width = 75.72255706787110
$sh17.Height = 22.251338582677167

# "ZoneTexte 140" - same shrink, text shortened to "reflexes" (position stays)
$sh25 = $s.Shapes.Item("ZoneTexte 140")
$sh25.TextFrame.TextRange.Text = "reflexes"
$sh25.Width = 75.72255706787110
$sh25.Height = 22.251338582677167
